$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find and delete rows for players removed from player_per_game_df:
# Brandon Ingram, Pascal Siakam, Zach Randolph
$playersToRemove = @("Brandon Ingram", "Pascal Siakam", "Zach Randolph")

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($i = $rowCount; $i -ge 1; $i--) {
    $cellValue = $ws.Cells.Item($i, 1).Value()
    if ($playersToRemove -contains $cellValue) {
        $ws.Rows.Item($i).Delete()
    }
}
